# Auto-generated edit script
# Applies cell-level numeric updates to the Kujata_Profits market-data sheets
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H32").Value = 1563
$ws.Range("J32").Value = 1789.625
$ws.Range("L32").Value = 1789.625
$ws.Range("N32").Value = -2441.625
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 25004000
$ws.Range("I113").Value = 40002800
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 40002800
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -39999546
$ws.Range("N113").Value = -12508
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
$ws.Range("H137").Value = 1435.8667
$ws.Range("I137").Value = 1129.7142
$ws.Range("J137").Value = 1703.75
$ws.Range("K137").Value = 3389.1426
$ws.Range("L137").Value = 5111.25
$ws.Range("M137").Value = -839.1425999999997
$ws.Range("N137").Value = -10211.25
$ws.Range("H138").Value = 2074.49
$ws.Range("I138").Value = 1484.3077
$ws.Range("J138").Value = 2162.6782
$ws.Range("K138").Value = 4452.9231
$ws.Range("L138").Value = 6488.034599999999
$ws.Range("M138").Value = 687.0769
$ws.Range("N138").Value = -16768.0346

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 1057.75
$ws.Range("I61").Value = 868.9231
$ws.Range("J61").Value = 1876
$ws.Range("K61").Value = 868.9231
$ws.Range("L61").Value = 1876
$ws.Range("M61").Value = -656.9231
$ws.Range("N61").Value = -2300
$ws.Range("H74").Value = 911.1579
$ws.Range("I74").Value = 893.625
$ws.Range("J74").Value = 1004.6667
$ws.Range("K74").Value = 893.625
$ws.Range("L74").Value = 1004.6667
$ws.Range("M74").Value = -19.625
$ws.Range("N74").Value = -2752.6667
$ws.Range("H77").Value = 911.1579
$ws.Range("I77").Value = 893.625
$ws.Range("J77").Value = 1004.6667
$ws.Range("K77").Value = 4468.125
$ws.Range("L77").Value = 5023.3335
$ws.Range("M77").Value = -100.125
$ws.Range("N77").Value = -13759.3335
$ws.Range("H115").Value = 37890
$ws.Range("J115").Value = 37890
$ws.Range("L115").Value = 37890
$ws.Range("N115").Value = -41024
$ws.Range("H123").Value = 68333.336
$ws.Range("J123").Value = 68333.336
$ws.Range("L123").Value = 68333.336
$ws.Range("N123").Value = -78133.336
$ws.Range("H126").Value = 9999.5
$ws.Range("I126").Value = 9999.5
$ws.Range("K126").Value = 29998.5
$ws.Range("M126").Value = -27528.5
$ws.Range("H136").Value = 1057.75
$ws.Range("I136").Value = 868.9231
$ws.Range("J136").Value = 1876
$ws.Range("K136").Value = 2606.7693
$ws.Range("L136").Value = 5628
$ws.Range("M136").Value = -56.76929999999993
$ws.Range("N136").Value = -10728

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 35715070
$ws.Range("I94").Value = 41667250
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 41667250
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -41666799
$ws.Range("N94").Value = -2902
$ws.Range("H99").Value = 31251346
$ws.Range("I99").Value = 33334642
$ws.Range("K99").Value = 33334642
$ws.Range("M99").Value = -33333144
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H134").Value = 4967.407
$ws.Range("I134").Value = 1117.7826
$ws.Range("J134").Value = 27102.75
$ws.Range("K134").Value = 3353.3478
$ws.Range("L134").Value = 81308.25
$ws.Range("M134").Value = -818.3478
$ws.Range("N134").Value = -86378.25

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 369.3846
$ws.Range("I22").Value = 288.57144
$ws.Range("J22").Value = 463.66666
$ws.Range("K22").Value = 288.57144
$ws.Range("L22").Value = 463.66666
$ws.Range("M22").Value = 61.42856
$ws.Range("N22").Value = -1163.66666
$ws.Range("H31").Value = 1516.6207
$ws.Range("I31").Value = 1194.3636
$ws.Range("J31").Value = 1713.5555
$ws.Range("K31").Value = 1194.3636
$ws.Range("L31").Value = 1713.5555
$ws.Range("M31").Value = -899.3635999999999
$ws.Range("N31").Value = -2303.5555
$ws.Range("H34").Value = 1516.6207
$ws.Range("I34").Value = 1194.3636
$ws.Range("J34").Value = 1713.5555
$ws.Range("K34").Value = 1194.3636
$ws.Range("L34").Value = 1713.5555
$ws.Range("M34").Value = -992.3635999999999
$ws.Range("N34").Value = -2117.5555
$ws.Range("H58").Value = 1052.2273
$ws.Range("I58").Value = 1085.8235
$ws.Range("K58").Value = 1085.8235
$ws.Range("M58").Value = -882.8235
$ws.Range("H109").Value = 9861.666999999999
$ws.Range("J109").Value = 9861.666999999999
$ws.Range("L109").Value = 9861.666999999999
$ws.Range("N109").Value = -11941.667
$ws.Range("H134").Value = 22224168
$ws.Range("I134").Value = 27779652
$ws.Range("J134").Value = 2233.3333
$ws.Range("K134").Value = 83338956
$ws.Range("L134").Value = 6699.999899999999
$ws.Range("M134").Value = -83336421
$ws.Range("N134").Value = -11769.9999
$ws.Range("H136").Value = 1052.2273
$ws.Range("I136").Value = 1085.8235
$ws.Range("K136").Value = 3257.4705
$ws.Range("M136").Value = -707.4704999999999

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H122").Value = 746
$ws.Range("I122").Value = 502.8
$ws.Range("J122").Value = 1093.4286
$ws.Range("K122").Value = 4525.2
$ws.Range("L122").Value = 9840.857399999999
$ws.Range("M122").Value = -2075.2
$ws.Range("N122").Value = -14740.8574
$ws.Range("H124").Value = 2933.3333
$ws.Range("I124").Value = 1100
$ws.Range("K124").Value = 3300
$ws.Range("M124").Value = 1610
$ws.Range("H126").Value = 5194.6
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -1060
$ws.Range("H131").Value = 24391752
$ws.Range("J131").Value = 1738.7742
$ws.Range("L131").Value = 5216.3226
$ws.Range("N131").Value = -15296.3226

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1361.8
$ws.Range("I16").Value = 997.5
$ws.Range("J16").Value = 2090.4
$ws.Range("K16").Value = 997.5
$ws.Range("L16").Value = 2090.4
$ws.Range("M16").Value = -827.5
$ws.Range("N16").Value = -2430.4
$ws.Range("H46").Value = 5420.7646
$ws.Range("J46").Value = 9866.888999999999
$ws.Range("L46").Value = 9866.888999999999
$ws.Range("N46").Value = -10242.889
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H107").Value = 66666
$ws.Range("I107").Value = 66666
$ws.Range("K107").Value = 66666
$ws.Range("M107").Value = -64746
$ws.Range("H136").Value = 1844.0769
$ws.Range("I136").Value = 1806.1818
$ws.Range("K136").Value = 5418.5454
$ws.Range("M136").Value = -2868.5454

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 3798.8572
$ws.Range("I132").Value = 3960.238
$ws.Range("K132").Value = 11880.714
$ws.Range("M132").Value = -9350.714
$ws.Range("H136").Value = 447.23254
$ws.Range("I136").Value = 363.17856
$ws.Range("J136").Value = 604.13336
$ws.Range("K136").Value = 1089.53568
$ws.Range("L136").Value = 1812.40008
$ws.Range("M136").Value = 1460.46432
$ws.Range("N136").Value = -6912.40008

Write-Output "Applied all cell updates."